$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values (before rotation) for rows 4, 5 and 6 across columns A..Q
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

$row4 = @{}
$row5 = @{}
$row6 = @{}

foreach ($col in $cols) {
    $row4[$col] = $ws.Range("${col}4").Value2
    $row5[$col] = $ws.Range("${col}5").Value2
    $row6[$col] = $ws.Range("${col}6").Value2
}

# Ensure columns are treated as plain text so values such as dates
# ("2022-11-01") are not auto-converted into Excel date serials.
foreach ($col in $cols) {
    $ws.Range("${col}4").NumberFormat = "@"
    $ws.Range("${col}5").NumberFormat = "@"
    $ws.Range("${col}6").NumberFormat = "@"
}

# Rotate rows: new row4 = old row6, new row5 = old row4, new row6 = old row5
foreach ($col in $cols) {
    $ws.Range("${col}4").Value = $row6[$col]
    $ws.Range("${col}5").Value = $row4[$col]
    $ws.Range("${col}6").Value = $row5[$col]
}
